$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. These cells already held text (not numeric)
# values such as "275.17" or "-1.41%" in the source workbook, so the
# NumberFormat must be forced to Text ("@") before assigning, otherwise
# Excel would reinterpret them as numbers/percentages.
$updates = @{
    "D2"  = "274.74";    "E2"  = "-1.70%"
    "D3"  = "27.44";     "E3"  = "1.62%"
    "D4"  = "4.803";     "E4"  = "-1.97%"
    "D5"  = "0.06291";   "E5"  = "-0.79%"
    "D7"  = "1.288";     "E7"  = "35.63%"
    "D8"  = "0.8753";    "E8"  = "-0.87%"
    "D9"  = "0.1525";    "E9"  = "3.73%"
    "D10" = "0.05023";   "E10" = "-3.15%"
    "D11" = "0.07485";   "E11" = "0.86%"
    "D12" = "0.02907";   "E12" = "-7.57%"
    "D13" = "0.09043";   "E13" = "-0.29%"
    "D14" = "0.001560";  "E14" = "0.22%"
    "D15" = "0.0006367"; "E15" = "1.27%"
    "D16" = "0.005923";  "E16" = "2.42%"
    "D17" = "3.450";     "E17" = "-0.80%"
    "D18" = "3.314";     "E18" = "-1.18%"
    "E19" = "-1.01%"
    "E20" = "0.60%"
    "D21" = "0.1318";    "E21" = "0.60%"
    "D22" = "3.900";     "E22" = "0.63%"
    "D23" = "0.04409";   "E23" = "2.09%"
    "D24" = "0.001170";  "E24" = "-0.65%"
    "D25" = "0.003836";  "E25" = "6.01%"
    "D26" = "0.0001199"; "E26" = "0.07%"
    "D27" = "0.0001936"; "E27" = "14.31%"
    "D40" = "0.04108";   "E40" = "1.47%"
    "D41" = "0.007040";  "E41" = "6.29%"
    "D42" = "0.1172";    "E42" = "0.54%"
    "D43" = "0.002019";  "E43" = "-13.97%"
    "D44" = "0.01122";   "E44" = "-9.67%"
    "D45" = "0.00005181";"E45" = "-0.49%"
    "D46" = "0.01999";   "E46" = "-11.11%"
    "E47" = "-37.51%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
